$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Produkty"

# Update headers
$ws.Range("B1").Value = "Produkt"
$ws.Range("D1").Value = "Cena Heureka (Kč)"
$ws.Range("E1").Value = "Moje cena"

# Row 2
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "JBL Tune 720BT"
$ws.Range("C2").Value = "https://sluchatka.heureka.cz/jbl-tune-720bt/#prehled/?sort-filter=lowest_price"
$ws.Range("D2").Value = 989
$ws.Range("E2").Value = 150

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "JBL Tune 720BT"
$ws.Range("C3").Value = "https://sluchatka.heureka.cz/jbl-tune-720bt/#prehled/?sort-filter=lowest_price"
$ws.Range("D3").Value = 989
$ws.Range("E3").Value = 150

# Row 4
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "JBL Tune 720BT"
$ws.Range("C4").Value = "https://sluchatka.heureka.cz/jbl-tune-720bt/#prehled/?sort-filter=lowest_price"
$ws.Range("D4").Value = 989
$ws.Range("E4").Value = 150

# Row 5
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "JBL Charge 6"
$ws.Range("C5").Value = "https://bluetooth-reproduktory.heureka.cz/jbl-charge-6/#prehled/?sort-filter=lowest_price"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 150

# Remove old column F which is no longer used
$ws.Range("F1:F5").Clear()
